# Append three new sentences/runs to the paragraph that currently ends with
# "...must be printed on the serial monitor/terminal. " in the
# "Inputs / Outputs" section.
#
# Target resulting runs (within the same paragraph):
#   "The module will also receive a series of image values from the object
#    detection module which will be given to "
#   "Matlab"                                   (wrapped by proofErr markers
#                                                in the reference document)
#   " to plot. "

$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute(
    "serial monitor/terminal. ", # FindText
    $true,                       # MatchCase
    $false,                      # MatchWholeWord
    $false,                      # MatchWildcards
    $false,                      # MatchSoundsLike
    $false,                      # MatchAllWordForms
    $true,                       # Forward
    1,                           # Wrap (wdFindContinue)
    $false,                      # Format
    "",                          # ReplaceWith
    0                            # Replace (wdReplaceNone)
)

if ($found) {
    # Move to the end of the matched text, then append the new runs one at
    # a time so each becomes its own <w:r> element (mirroring the diff).
    $rng.Collapse(0)
    $rng.InsertAfter("The module will also receive a series of image values from the object detection module which will be given to ")

    $rng.Collapse(0)
    $rng.InsertAfter("Matlab")

    $rng.Collapse(0)
    $rng.InsertAfter(" to plot. ")
}
